$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.971.53"
$ws.Range("E2").Value = "  +5.16%  "
$ws.Range("D3").Value = "2.337.98"
$ws.Range("E3").Value = "  +4.58%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").Value = "'305.82"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'96.63"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +6.56%  "
$ws.Range("D10").Value = "'35.76"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("D12").Value = "'7.44"
$ws.Range("E12").Value = "  +5.69%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "2.691.52"
$ws.Range("E14").Value = "  +4.35%  "
$ws.Range("D15").Value = "2.336.63"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'14.16"
$ws.Range("E16").Value = "  +6.05%  "
$ws.Range("D17").Value = "'0.834"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "46.806.38"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  +18.40%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "'67.67"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").Value = "'254.15"
$ws.Range("E23").Value = "  +7.82%  "
$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "'42.27"
$ws.Range("E27").Value = "  +15.18%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").Value = "'20.30"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "'5.81"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'0.0817"
$ws.Range("E32").Value = "  +6.94%  "
$ws.Range("D33").Value = "'146.48"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").Value = "'2.61"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'3.16"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "'1.82"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'4.00"
$ws.Range("E39").Value = "  +8.22%  "
$ws.Range("D40").Value = "'0.0313"
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").Value = "'13.89"
$ws.Range("E42").Value = "  -5.80%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  +13.29%  "
$ws.Range("D45").Value = "1.817.40"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'91.85"
$ws.Range("E46").Value = "  +14.26%  "
$ws.Range("D47").Value = "'74.95"
$ws.Range("E47").Value = "  +10.46%  "
$ws.Range("E48").Value = "  +6.11%  "
$ws.Range("D49").Value = "'99.01"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'4.84"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'55.12"
$ws.Range("E51").Value = "  +3.67%  "
